$d = $word.ActiveDocument

# Helper: find the paragraph whose text contains a given marker substring.
function Get-ParaByMarker($marker) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text -like "*$marker*") {
            return $p
        }
    }
    return $null
}

# --- 1) "Serial.printF" paragraph: wrap the run in spellStart/spellEnd proofErr ---
$p1 = Get-ParaByMarker("Serial.printF")
if ($null -eq $p1) { throw "paragraph with 'Serial.printF' not found" }
$xml1 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:pPr><w:ind w:left="687"/></w:pPr>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Serial.printF</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t xml:space="preserve"> ((&#8220;Current Temp is, %i /n&#8221;), currentTemp) { </w:t></w:r>' +
        '</w:p>'
[void]$p1.Range.InsertXML($xml1)

# --- 2) "bandpassing" paragraph: split run, wrap "bandpassing" in spellStart/spellEnd ---
$p2 = Get-ParaByMarker("bandpassing")
if ($null -eq $p2) { throw "paragraph with 'bandpassing' not found" }
$xml2 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:pPr><w:ind w:left="2160" w:firstLine="5"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
        '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">//angulate &amp; tune: follow noise to source; angle drone //listening by </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>bandpassing</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> - once tuned in focus cameras</w:t></w:r>' +
        '</w:p>'
[void]$p2.Range.InsertXML($xml2)

# --- 3) "If (no_anomaly){" paragraph: split run, add spellStart/spellEnd + gramStart/gramEnd ---
$p3 = Get-ParaByMarker("no_anomaly")
if ($null -eq $p3) { throw "paragraph with 'no_anomaly' not found" }
$xml3 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:pPr><w:ind w:left="687"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
        '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>If (</w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>no_</w:t></w:r>' +
        '<w:proofErr w:type="gramStart"/>' +
        '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>anomaly</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>){</w:t></w:r>' +
        '<w:proofErr w:type="gramEnd"/>' +
        '</w:p>'
[void]$p3.Range.InsertXML($xml3)

# --- 4) "buttonPush_3_long" paragraph: rename to "longPush", drop trailing space run,
#         add spellStart/spellEnd proofErr and en-US lang on both pPr and run ---
$p4 = Get-ParaByMarker("buttonPush_3_long")
if ($null -eq $p4) { throw "paragraph with 'buttonPush_3_long' not found" }
$xml4 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:pPr><w:ind w:left="730" w:firstLine="710"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:rPr><w:u w:val="single" w:color="000000"/><w:lang w:val="en-US"/></w:rPr><w:t>longPush</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '</w:p>'
[void]$p4.Range.InsertXML($xml4)
